$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("max") entirely; this shifts D->C and E->D.
$ws.Range("C1").EntireColumn.Delete()

# Update the values that differ from a simple shift.
$ws.Range("B2").Value = 8.03068275202429
$ws.Range("C2").Value = "s__CADAUA01 sp900315375"
